$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 30367.75
$ws.Range("J114").Value = 30367.75
$ws.Range("L114").Value = 30367.75
$ws.Range("N114").Value = -39045.75
$ws.Range("H117").Value = 48734
$ws.Range("J117").Value = 48734
$ws.Range("L117").Value = 48734
$ws.Range("N117").Value = -57912
$ws.Range("H124").Value = 55972.2
$ws.Range("J124").Value = 55972.2
$ws.Range("L124").Value = 55972.2
$ws.Range("N124").Value = -65792.2
$ws.Range("H128").Value = 54985
$ws.Range("J128").Value = 54985
$ws.Range("L128").Value = 54985
$ws.Range("N128").Value = -64945
$ws.Range("H132").Value = 21293.104
$ws.Range("I132").Value = 3378.8838
$ws.Range("K132").Value = 10136.6514
$ws.Range("M132").Value = -7606.651400000001
$ws.Range("H138").Value = 2587.898
$ws.Range("I138").Value = 2032.4062
$ws.Range("J138").Value = 2857.2273
$ws.Range("K138").Value = 6097.2186
$ws.Range("L138").Value = 8571.6819
$ws.Range("M138").Value = -957.2186000000002
$ws.Range("N138").Value = -18851.6819
$ws.Range("H141").Value = 3592.4546
$ws.Range("I141").Value = 2317.7222
$ws.Range("K141").Value = 6953.1666
$ws.Range("M141").Value = -1773.1666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 44655
$ws.Range("J111").Value = 44655
$ws.Range("L111").Value = 44655
$ws.Range("N111").Value = -52835
$ws.Range("H113").Value = 42195.6
$ws.Range("J113").Value = 42195.6
$ws.Range("L113").Value = 42195.6
$ws.Range("N113").Value = -50873.6
$ws.Range("H114").Value = 30313.5
$ws.Range("J114").Value = 30313.5
$ws.Range("L114").Value = 30313.5
$ws.Range("N114").Value = -38991.5
$ws.Range("H117").Value = 44273.168
$ws.Range("J117").Value = 44273.168
$ws.Range("L117").Value = 44273.168
$ws.Range("N117").Value = -53451.168
$ws.Range("H118").Value = 49998
$ws.Range("J118").Value = 49998
$ws.Range("L118").Value = 49998
$ws.Range("N118").Value = -53312
$ws.Range("H119").Value = 51359.332
$ws.Range("J119").Value = 51359.332
$ws.Range("L119").Value = 51359.332
$ws.Range("N119").Value = -61035.332
$ws.Range("H123").Value = 51425
$ws.Range("J123").Value = 51425
$ws.Range("L123").Value = 51425
$ws.Range("N123").Value = -61225
$ws.Range("H125").Value = 39715
$ws.Range("J125").Value = 39715
$ws.Range("L125").Value = 39715
$ws.Range("N125").Value = -49555
$ws.Range("H130").Value = 44873.668
$ws.Range("J130").Value = 44873.668
$ws.Range("L130").Value = 44873.668
$ws.Range("N130").Value = -54913.668
$ws.Range("H131").Value = 49907
$ws.Range("J131").Value = 49907
$ws.Range("L131").Value = 49907
$ws.Range("N131").Value = -59987

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 46663.25
$ws.Range("J108").Value = 46663.25
$ws.Range("L108").Value = 46663.25
$ws.Range("N108").Value = -54343.25
$ws.Range("H110").Value = 47992
$ws.Range("J110").Value = 47992
$ws.Range("L110").Value = 47992
$ws.Range("N110").Value = -56172
$ws.Range("H111").Value = 33567.332
$ws.Range("J111").Value = 33567.332
$ws.Range("L111").Value = 33567.332
$ws.Range("N111").Value = -41747.332
$ws.Range("H112").Value = 45153.668
$ws.Range("J112").Value = 45153.668
$ws.Range("L112").Value = 45153.668
$ws.Range("N112").Value = -48107.668
$ws.Range("H117").Value = 49914
$ws.Range("J117").Value = 49914
$ws.Range("L117").Value = 49914
$ws.Range("N117").Value = -59092
$ws.Range("H124").Value = 44735.75
$ws.Range("J124").Value = 44735.75
$ws.Range("L124").Value = 44735.75
$ws.Range("N124").Value = -54555.75
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H126").Value = 50772
$ws.Range("J126").Value = 50772
$ws.Range("L126").Value = 50772
$ws.Range("N126").Value = -60652
$ws.Range("H130").Value = 53387.5
$ws.Range("J130").Value = 53387.5
$ws.Range("L130").Value = 53387.5
$ws.Range("N130").Value = -63427.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49839.668
$ws.Range("J20").Value = 49839.668
$ws.Range("L20").Value = 49839.668
$ws.Range("N20").Value = -50311.668
$ws.Range("H30").Value = 49839.668
$ws.Range("J30").Value = 49839.668
$ws.Range("L30").Value = 49839.668
$ws.Range("N30").Value = -50021.668
$ws.Range("H110").Value = 42999
$ws.Range("J110").Value = 42999
$ws.Range("L110").Value = 42999
$ws.Range("N110").Value = -51179
$ws.Range("H112").Value = 40492
$ws.Range("J112").Value = 40492
$ws.Range("L112").Value = 40492
$ws.Range("N112").Value = -43446
$ws.Range("H116").Value = 47313.332
$ws.Range("J116").Value = 47313.332
$ws.Range("L116").Value = 47313.332
$ws.Range("N116").Value = -56491.332
$ws.Range("H119").Value = 48753
$ws.Range("J119").Value = 48753
$ws.Range("L119").Value = 48753
$ws.Range("N119").Value = -58429
$ws.Range("H128").Value = 49839.668
$ws.Range("J128").Value = 49839.668
$ws.Range("L128").Value = 49839.668
$ws.Range("N128").Value = -59799.668

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 17363.084
$ws.Range("I107").Value = 13169.125
$ws.Range("J107").Value = 25751
$ws.Range("K107").Value = 39507.375
$ws.Range("L107").Value = 77253
$ws.Range("M107").Value = -37587.375
$ws.Range("N107").Value = -81093
$ws.Range("H131").Value = 3949.162
$ws.Range("I131").Value = 13091.125
$ws.Range("J131").Value = 1427.2413
$ws.Range("K131").Value = 39273.375
$ws.Range("L131").Value = 4281.7239
$ws.Range("M131").Value = -34233.375
$ws.Range("N131").Value = -14361.7239

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 37311.668
$ws.Range("J114").Value = 37311.668
$ws.Range("L114").Value = 37311.668
$ws.Range("N114").Value = -45989.668
$ws.Range("H116").Value = 49734
$ws.Range("J116").Value = 49734
$ws.Range("L116").Value = 49734
$ws.Range("N116").Value = -58912
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H114").Value = 25920.4
$ws.Range("J114").Value = 25920.4
$ws.Range("L114").Value = 25920.4
$ws.Range("N114").Value = -34598.4
$ws.Range("H116").Value = 51680
$ws.Range("J116").Value = 51680
$ws.Range("L116").Value = 51680
$ws.Range("N116").Value = -60858
$ws.Range("H118").Value = 43409
$ws.Range("J118").Value = 43409
$ws.Range("L118").Value = 43409
$ws.Range("N118").Value = -46723
$ws.Range("H120").Value = 57491
$ws.Range("J120").Value = 57491
$ws.Range("L120").Value = 57491
$ws.Range("N120").Value = -67167
$ws.Range("H124").Value = 42210
$ws.Range("J124").Value = 42210
$ws.Range("L124").Value = 42210
$ws.Range("N124").Value = -52030
$ws.Range("H125").Value = 40211
$ws.Range("J125").Value = 40211
$ws.Range("L125").Value = 40211
$ws.Range("N125").Value = -50051
$ws.Range("H128").Value = 37689.5
$ws.Range("J128").Value = 37689.5
$ws.Range("L128").Value = 37689.5
$ws.Range("N128").Value = -47649.5
$ws.Range("H130").Value = 47570
$ws.Range("J130").Value = 47570
$ws.Range("L130").Value = 47570
$ws.Range("N130").Value = -57610

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H110").Value = 48883.5
$ws.Range("J110").Value = 48883.5
$ws.Range("L110").Value = 48883.5
$ws.Range("N110").Value = -57063.5
$ws.Range("H116").Value = 48690.668
$ws.Range("J116").Value = 48690.668
$ws.Range("L116").Value = 48690.668
$ws.Range("N116").Value = -57868.668
$ws.Range("H117").Value = 49301
$ws.Range("J117").Value = 49301
$ws.Range("L117").Value = 49301
$ws.Range("N117").Value = -58479
$ws.Range("H121").Value = 43589.5
$ws.Range("J121").Value = 43589.5
$ws.Range("L121").Value = 43589.5
$ws.Range("N121").Value = -47083.5
$ws.Range("H132").Value = 1264.4182
$ws.Range("I132").Value = 1056.9318
$ws.Range("J132").Value = 2094.3635
$ws.Range("K132").Value = 3170.7954
$ws.Range("L132").Value = 6283.0905
$ws.Range("M132").Value = -640.7954
$ws.Range("N132").Value = -11343.0905
$ws.Range("N108").ClearContents()
